# Feigenbaum & Sensitivität - minor adjustments
#
# 1) On the "Inhalt" overview slide (slide 2), add a new sub-bullet
#    "Grundlagen" (indent level 3 / OOXML lvl="2") right after the
#    "Feigenbaum-Diagramm" bullet, before "Ordnung".
# 2) On the "Sensitivität und Chaos" slide (slide 7), fix a wording typo:
#    "Haben kleine Änderungen..." -> "Haben sehr kleine Änderungen..."

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Slide 2 - "Inhalt" content placeholder: insert "Grundlagen" bullet
# ---------------------------------------------------------------------
$s1 = $p.Slides.Item(2)
$sh1 = $s1.Shapes.Item(2)
$tr1 = $sh1.TextFrame.TextRange

# Rebuild the bullet list text with the new "Grundlagen" line inserted
# right after "Feigenbaum-Diagramm". (Re-assigning .Text is the reliable
# way to introduce a genuinely new paragraph in this host; InsertBefore/
# InsertAfter merely splice a literal line-break into the existing run.)
$tr1.Text = "Feigenbaum`rFeigenbaum-Diagramm`rGrundlagen`rOrdnung`rChaos`rFeigenbaum-Konstante`t`rSensitivität`rArten chaotischen Verhaltens`rSchmetterlingseffekt`t"

# Restore the outline indent levels for every paragraph (re-assigning
# .Text above resets them all back to the top level).
$full1 = $sh1.TextFrame.TextRange.Text

function Set-Indent($textRange, $fullText, $word, $level) {
    $idx = $fullText.IndexOf($word)
    $run = $textRange.Characters($idx + 1, $word.Length)
    $run.IndentLevel = $level
}

Set-Indent $tr1 $full1 "Feigenbaum-Diagramm" 2
Set-Indent $tr1 $full1 "Grundlagen" 3
Set-Indent $tr1 $full1 "Ordnung" 3
Set-Indent $tr1 $full1 "Chaos" 3
Set-Indent $tr1 $full1 "Feigenbaum-Konstante" 2
Set-Indent $tr1 $full1 "Arten chaotischen Verhaltens" 2
Set-Indent $tr1 $full1 "Schmetterlingseffekt" 2

# ---------------------------------------------------------------------
# 2) Slide 7 - "Sensitivität und Chaos": tweak wording
# ---------------------------------------------------------------------
$s2 = $p.Slides.Item(7)
$sh2 = $s2.Shapes.Item(2)
$tr2 = $sh2.TextFrame.TextRange

$oldText = "Haben kleine Änderungen stark unterschiedliche Ergebnisse zu folge wirkt das System chaotisch"
$newText = "Haben sehr kleine Änderungen stark unterschiedliche Ergebnisse zu folge wirkt das System chaotisch"

$full2 = $tr2.Text
$pos = $full2.IndexOf($oldText)
if ($pos -ge 0) {
    $target = $tr2.Characters($pos + 1, $oldText.Length)
    $target.Text = $newText
}
